# completed part 3 of design details
#
# Fills in the five "label:" paragraphs under heading 3 ("API:", "Scraper:",
# "Database:", "Frontend:", "Deployment:") with their write-up text, bolding
# each label. The "_GoBack" bookmark (which Word maintains at the location of
# the most recent edit) moves from just before "Scraper:" to inside the new
# "Deployment:" text, right after the word "running".

$d = $word.ActiveDocument

$wNs = 'xmlns:w="http://schemas.openxmlformats.org/wordprocessingml/2006/main"'

# Builds a <w:r> run. $Bold=$true adds <w:b/><w:bCs/>; $Preserve=$true adds
# xml:space="preserve". (Plain [bool] params are used instead of [switch] -
# this host's PowerShell-subset doesn't reliably evaluate switch parameters.)
function New-RunXml([string]$Text, [bool]$Bold = $false, [bool]$Preserve = $false) {
    $rPr = ""
    if ($Bold) { $rPr = "<w:rPr><w:b/><w:bCs/></w:rPr>" }
    $space = ""
    if ($Preserve) { $space = ' xml:space="preserve"' }
    $escaped = $Text.Replace("&", "&amp;").Replace("<", "&lt;").Replace(">", "&gt;")
    return "<w:r>$rPr<w:t$space>$escaped</w:t></w:r>"
}

# Replaces the whole paragraph at 1-based Paragraphs index $ParaIndex with a
# paragraph carrying the same "no extra spacing" formatting and the supplied
# run/bookmark markup, via Range.InsertXML (a WordprocessingML import, same
# mechanism Word itself uses for "Keep Source Formatting" pastes).
function Set-ParagraphRuns([int]$ParaIndex, [string]$InnerXml) {
    $rng = $d.Paragraphs.Item($ParaIndex).Range
    $pPr = "<w:pPr><w:spacing w:after=`"0`" w:afterAutospacing=`"0`"/></w:pPr>"
    $pkg = '<pkg:package xmlns:pkg="http://schemas.microsoft.com/office/2006/xmlPackage">' +
           '<pkg:part pkg:name="/word/document.xml" pkg:contentType="application/vnd.openxmlformats-officedocument.wordprocessingml.document.main+xml">' +
           "<pkg:xmlData><w:document $wNs><w:body><w:p>$pPr$InnerXml</w:p></w:body></w:document></pkg:xmlData>" +
           '</pkg:part></pkg:package>'
    $rng.InsertXML($pkg)
}

# Locate the five label paragraphs by their current text so this keeps
# working even if unrelated paragraphs are added/removed earlier in the doc.
# Paragraph.Range.Text includes the trailing paragraph mark, so trim it
# before comparing.
function Find-ParagraphIndex([string]$LabelText) {
    for ($i = 1; $i -le $d.Paragraphs.Count; $i++) {
        $t = $d.Paragraphs.Item($i).Range.Text.TrimEnd([char]13, [char]7)
        if ($t -eq $LabelText) {
            return $i
        }
    }
    return -1
}

$apiIdx = Find-ParagraphIndex "API:"
$scraperIdx = Find-ParagraphIndex "Scraper:"
$databaseIdx = Find-ParagraphIndex "Database:"
$frontendIdx = Find-ParagraphIndex "Frontend:"
$deploymentIdx = Find-ParagraphIndex "Deployment:"

foreach ($pair in @(@("API:", $apiIdx), @("Scraper:", $scraperIdx), @("Database:", $databaseIdx), @("Frontend:", $frontendIdx), @("Deployment:", $deploymentIdx))) {
    if ($pair[1] -lt 1) {
        throw ("Could not locate paragraph for label: " + $pair[0])
    }
}

# API:
$apiRuns = (New-RunXml "API:" -Bold) +
           (New-RunXml " Our API will be written in node JS" -Preserve)
Set-ParagraphRuns $apiIdx $apiRuns

# Scraper: (this also removes the old "_GoBack" bookmark that used to sit
# here, since the whole paragraph content is replaced)
$scraperRuns = (New-RunXml "Scraper:" -Bold) +
               (New-RunXml " Our web scraper will be built in node JS using the cheerio library. We chose this as the library is easy to use, relatively light weight, and it makes it extremely simple " -Preserve) +
               (New-RunXml "to create JSON objects from the scraped reports from a JavaScript program.") +
               (New-RunXml " This means that it will work seamlessly with the API and database." -Preserve)
Set-ParagraphRuns $scraperIdx $scraperRuns

# Database:
$databaseRuns = (New-RunXml "Database:" -Bold) +
                (New-RunXml " Our database will be constructed using SQLite3" -Preserve)
Set-ParagraphRuns $databaseIdx $databaseRuns

# Frontend:
$frontendRuns = (New-RunXml "Frontend:" -Bold) +
                (New-RunXml " " -Preserve) +
                (New-RunXml "The frontend for our web") +
                (New-RunXml " " -Preserve) +
                (New-RunXml "app will be constructed using typical web programming languages such as HTML5, CSS and " -Preserve) +
                (New-RunXml "JavaScript.")
Set-ParagraphRuns $frontendIdx $frontendRuns

# Deployment: - the new "_GoBack" bookmark lands right after "running"
$deploymentRuns = (New-RunXml "Deployment:" -Bold) +
                  (New-RunXml " Our API will be deployed " -Preserve) +
                  (New-RunXml "on the Heroku web") +
                  (New-RunXml " hosting service. This allows us to" -Preserve) +
                  (New-RunXml " easily" -Preserve) +
                  (New-RunXml " host our API on the internet, and run our scraper on a schedule for free" -Preserve) +
                  (New-RunXml ".") +
                  (New-RunXml " For local development, the API and the scraper will both be " -Preserve) +
                  (New-RunXml "running") +
                  '<w:bookmarkStart w:id="0" w:name="_GoBack"/><w:bookmarkEnd w:id="0"/>' +
                  (New-RunXml " on Linux operating systems." -Preserve)
Set-ParagraphRuns $deploymentIdx $deploymentRuns

Write-Output "Done: populated API/Scraper/Database/Frontend/Deployment paragraphs."
